$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44600
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 12500
$ws.Range("P2").Value = 208

# Row 3
$ws.Range("D3").Value = 44216
$ws.Range("J3").Value = 55
$ws.Range("M3").Value = 9773
$ws.Range("P3").Value = 163

# Row 4
$ws.Range("D4").Value = 44259
$ws.Range("J4").Value = 70
$ws.Range("M4").Value = 9214
$ws.Range("O4").Value = "Región del Maule"

# Row 5
$ws.Range("D5").Value = 44208
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 7350
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 122

# Row 6
$ws.Range("D6").Value = 44610
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11500
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 192

# Row 7
$ws.Range("D7").Value = 44253
$ws.Range("J7").Value = 95
$ws.Range("K7").Value = 9500
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9658
$ws.Range("P7").Value = 161

# Row 8
$ws.Range("D8").Value = 44224
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 8500
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8719
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 145

# Row 9
$ws.Range("D9").Value = 44615
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11500
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 192

# Row 10
$ws.Range("D10").Value = 44159
$ws.Range("J10").Value = 35
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7714
$ws.Range("P10").Value = 129

# Row 11
$ws.Range("D11").Value = 44687
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9500
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 158

# Row 12
$ws.Range("D12").Value = 44627
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 9500
$ws.Range("M12").Value = 9250
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 154

# Row 13
$ws.Range("D13").Value = 44608
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 208

# Row 15
$ws.Range("D15").Value = 44624
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = 10500
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 175

# Row 16
$ws.Range("D16").Value = 44690
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 167

# Row 17
$ws.Range("D17").Value = 44692
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 167

# Row 18
$ws.Range("D18").Value = 44596
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = 12500
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 208

# Row 19
$ws.Range("D19").Value = 44594
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12500
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 208

# Row 20
$ws.Range("D20").Value = 44204
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = 9500
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9722
$ws.Range("P20").Value = 162

# Row 21
$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 90
$ws.Range("K21").Value = 7500
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7667
$ws.Range("P21").Value = 128

# Row 22
$ws.Range("D22").Value = 44202
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8400
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 140

# Row 23
$ws.Range("D23").Value = 44162
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8500
$ws.Range("M23").Value = 8209
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 137

# Row 24
$ws.Range("D24").Value = 44630
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9500
$ws.Range("M24").Value = 9250
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 154

# Row 25
$ws.Range("D25").Value = 44264
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = 8500
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8709
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 145

# Row 26
$ws.Range("D26").Value = 44210
$ws.Range("J26").Value = 60
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 8417
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 140

# Row 27
$ws.Range("D27").Value = 44218
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 9615
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 160

# Row 28
$ws.Range("D28").Value = 44671
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 6500
$ws.Range("P28").Value = 108

# Row 29
$ws.Range("D29").Value = 44271
$ws.Range("J29").Value = 55
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 9500
$ws.Range("M29").Value = 9227
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 154

# New row 30
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44617
$ws.Range("D30").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112001
$ws.Range("G30").Value = "Berenjena"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 11000
$ws.Range("M30").Value = 10500
$ws.Range("N30").Value = "$/caja 60 unidades"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 175
$ws.Range("Q30").Value = 60
$ws.Range("R30").Value = "Hortaliza"
